$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.611.74'
$ws.Range("E2").Value = '  +10.51%  '
$ws.Range("D3").Value = '3.265.47'
$ws.Range("E3").Value = '  +6.66%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.12'
$ws.Range("E5").Value = '  +3.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '112.02'
$ws.Range("E6").Value = '  +9.82%  '
$ws.Range("E7").Value = '  +4.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").Value = '  +7.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.53'
$ws.Range("E10").Value = '  +7.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0961'
$ws.Range("E11").Value = '  +13.33%  '
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("D13").Value = '3.772.89'
$ws.Range("E13").Value = '  +6.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.34'
$ws.Range("E14").Value = '  +5.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.14'
$ws.Range("E15").Value = '  +6.05%  '
$ws.Range("D16").Value = '3.258.97'
$ws.Range("E16").Value = '  +6.54%  '
$ws.Range("E17").Value = '  +6.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.04'
$ws.Range("E18").Value = '  +4.13%  '
$ws.Range("D19").Value = '56.516.95'
$ws.Range("E19").Value = '  +10.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.34'
$ws.Range("E20").Value = '  +4.60%  '
$ws.Range("E21").Value = '  +10.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.13'
$ws.Range("E22").Value = '  +7.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '299.46'
$ws.Range("E23").Value = '  +13.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.30'
$ws.Range("E24").Value = '  +8.11%  '
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.15'
$ws.Range("E26").Value = '  +3.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.39'
$ws.Range("E27").Value = '  +5.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.37'
$ws.Range("E28").Value = '  +4.57%  '
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E30").Value = '  +4.42%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +7.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.11'
$ws.Range("E33").Value = '  +6.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.34'
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.13'
$ws.Range("E36").Value = '  +2.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.79'
$ws.Range("E37").Value = '  +3.57%  '
$ws.Range("E38").Value = '  +26.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.54'
$ws.Range("E39").Value = '  +6.01%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.72'
$ws.Range("E41").Value = '  +7.56%  '
$ws.Range("E42").Value = '  +6.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '132.70'
$ws.Range("E43").Value = '  +2.11%  '
$ws.Range("E44").Value = '  +4.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.00'
$ws.Range("E45").Value = '  +6.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.286'
$ws.Range("E46").Value = '  -2.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.26'
$ws.Range("E47").Value = '  +2.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.20'
$ws.Range("E48").Value = '  +56.45%  '
$ws.Range("D49").Value = '2.152.35'
$ws.Range("E49").Value = '  +4.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.10'
$ws.Range("E50").Value = '  +1.72%  '
$ws.Range("E51").Value = '  -3.55%  '
